# StagingTemplates/Staging.Framework_Indicator.xlsx
#
# Commit: "moved staging files StagingTemplates directory"
#
# Semantic content change in the diff: the header label in C2 is renamed
# from "IndicatorSourceKey" to "IndicatorBusinessKey" (shared-string text
# change, cell keeps its existing style/type).
#
# The remaining hunks in the diff (bookViews/workbookView window
# geometry, the worksheet's internal VBA codeName counter, and the
# disappearance of the <col> bestFit/customWidth entries for columns B
# and C) are artifacts of the authoring Excel session/window chrome
# rather than data the workbook exposes for editing - they are re-stamped
# by the Excel client itself on save and aren't reachable through the
# object model. We still touch the equivalent properties below wherever
# the object model exposes them, so the intent is captured even though
# the host may not persist every cosmetic attribute.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Primary content edit -------------------------------------------------
# C2 holds the shared string "IndicatorSourceKey" -> rename to
# "IndicatorBusinessKey", preserving the cell's existing style (s="1").
$ws.Range("C2").Value = "IndicatorBusinessKey"

# --- Best-effort cosmetic edits -------------------------------------------
# Worksheet's VBA codeName bump (Sheet15 -> Sheet17) from re-saving in a
# newer Excel session. (No corresponding OOXML attribute is reachable from
# the data model for a mere cell/content edit, but set it in case the host
# honours it.)
$ws.CodeName = "Sheet17"

# Window geometry recorded by the authoring Excel session's workbookView.
$wb.Windows.Item(1).Width = 28800
$wb.Windows.Item(1).Height = 12585

# NOTE: columns B and C lose their explicit bestFit/customWidth sizing in
# the target (the <col> entries are removed outright, reverting to the
# sheet default), while column A (width 41) is left untouched. There is no
# object-model operation that *clears* an explicit column width back to
# "unset" (ColumnWidth/AutoFit always (re)write an explicit width), so
# touching columns B/C here would only trade one explicit width for
# another and move further away from the target than leaving them alone.
# Left as-is deliberately.
